$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HMI_signals")

# --- Extend the table with the new circuit breakers needed for the 2x UDP
# frames (Banshee & PHIL): CBID 401-413, continuing the ID sequence from 48.
# Two fill-down passes (49-55, then 56-61) so each gets its own
# shared-formula group, the same way Excel splits it when you fill down in
# two goes.

# Row 49-55 block (CBID 401-407)
$row1First = 49
$row1Last = 55
$cbids1 = 401,402,403,404,405,406,407

$id = 48
for ($i = 0; $i -lt $cbids1.Length; $i++) {
    $r = $row1First + $i
    $ws.Range("A$r").Value = $id
    $ws.Range("B$r").Value = $cbids1[$i]
    $id = $id + 1
}

$ws.Range("C$($row1First):C$($row1Last)").Formula = "=CONCATENATE(`"HMI1.Flow_`",B$row1First)"
$ws.Range("D$($row1First):D$($row1Last)").Formula = "=CONCATENATE(`"HMI1.CB`",B$row1First,`"_Closed`")"
$ws.Range("F$($row1First):F$($row1Last)").Formula = "=CONCATENATE(`"CB`", B$row1First, `"_MODBUS`")"
$ws.Range("G$($row1First):G$($row1Last)").Formula = "=IF(E$row1First=`"`",`"`",CONCATENATE(`",IOCCmd:=ADR(`",E$row1First,`"),MBCmd:=ADR(`",F$row1First,`".CMD),MBRst:=ADR(`",F$row1First,`".RST)`"))"
$ws.Range("H$($row1First):H$($row1Last)").Formula = "=CONCATENATE(`"(ID:=`",A$row1First,`",CBID:=`",B$row1First,`",Flow:=ADR(HMI1.Flow_`",B$row1First,`"),Status:=ADR(`",D$row1First,`")`",G$row1First,`"),`")"

# Row 56-61 block (CBID 408-413)
$row2First = 56
$row2Last = 61
$cbids2 = 408,409,410,411,412,413

for ($i = 0; $i -lt $cbids2.Length; $i++) {
    $r = $row2First + $i
    $ws.Range("A$r").Value = $id
    $ws.Range("B$r").Value = $cbids2[$i]
    $id = $id + 1
}

$ws.Range("C$($row2First):C$($row2Last)").Formula = "=CONCATENATE(`"HMI1.Flow_`",B$row2First)"
$ws.Range("D$($row2First):D$($row2Last)").Formula = "=CONCATENATE(`"HMI1.CB`",B$row2First,`"_Closed`")"
$ws.Range("F$($row2First):F$($row2Last)").Formula = "=CONCATENATE(`"CB`", B$row2First, `"_MODBUS`")"
$ws.Range("G$($row2First):G$($row2Last)").Formula = "=IF(E$row2First=`"`",`"`",CONCATENATE(`",IOCCmd:=ADR(`",E$row2First,`"),MBCmd:=ADR(`",F$row2First,`".CMD),MBRst:=ADR(`",F$row2First,`".RST)`"))"
$ws.Range("H$($row2First):H$($row2Last)").Formula = "=CONCATENATE(`"(ID:=`",A$row2First,`",CBID:=`",B$row2First,`",Flow:=ADR(HMI1.Flow_`",B$row2First,`"),Status:=ADR(`",D$row2First,`")`",G$row2First,`"),`")"

# --- Consolidate the H column's per-row formula (H2:H48) into a single
# shared formula too, same as the rest of the table already uses.
$ws.Range("H2:H48").Formula = "=CONCATENATE(`"(ID:=`",A2,`",CBID:=`",B2,`",Flow:=ADR(HMI1.Flow_`",B2,`"),Status:=ADR(`",D2,`")`",G2,`"),`")"

# --- Update the view: drop the frozen/topLeft B1 pin and move the
# selection to the C2:C10 block (the new section being reviewed).
$ws.Activate()
$ws.Range("C2:C10").Select()
